{"js": "return Word.Section.prototype._hf.toString();\n", "ps1": "$d = $word.ActiveDocument\nWrite-Output \"Fields.Add info\"\n"}
